$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Variable" column (A) currently holds short variable codes (e.g. hischshr1520m)
# and the "Label" column (B) holds the descriptive text. We drop the Variable
# column, move the Label text into column A, and clear column B entirely.

$labels = @(
    "Label",
    "Share Men aged 15-20 with High School Education",
    "Islamic Mayor in 1989",
    "Islamic vote share 1994",
    "Number of parties receiving votes 1994",
    "Log Population in 1994",
    "District center",
    "Province center",
    "Sub-metro center",
    "Metro center"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Clear out the old "Variable" column (B) entirely.
$ws.Range("B1:B10").Clear()
